# Fixing the big mistake: correct erroneous Total/Community/IGA consumption
# values (columns B, D, E) for rows 2-13 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = @{ B = 11308.81935935001;  D = 606.2584714166666;  E = 2028.987443183333 }
    3  = @{ B = 10675.96777400001;  D = 572.2136810166667;  E = 1825.988489983333 }
    4  = @{ B = 11456.17485328334;  D = 625.20078605;        E = 2040.56723405 }
    5  = @{ B = 11057.11805373335;  D = 599.6472717;         E = 2112.526876 }
    6  = @{ B = 11231.07623018335;  D = 609.2900709666667;  E = 2066.675799533333 }
    7  = @{ B = 11017.75471855001;  D = 588.2863521833333;  E = 2116.900388666667 }
    8  = @{ B = 11343.05436751668;  D = 627.9682417833333;  E = 2132.2404149 }
    9  = @{ B = 11467.18355501668;  D = 624.0573869999999;  E = 2187.1635299 }
    10 = @{ B = 10974.30720706668;  D = 609.8610163833334;  E = 1906.849093216667 }
    11 = @{ B = 11438.48934851668;  D = 632.4839693833334;  E = 2131.58672435 }
    12 = @{ B = 11036.06676451668;  D = 596.8665514666667;  E = 2073.6274182 }
    13 = @{ B = 11090.75984615001;  D = 610.0221540333333;  E = 2129.783974116667 }
}

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    $ws.Range("B$row").Value = $vals.B
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("E$row").Value = $vals.E
}
